# edit.ps1 - applies the "finish input and output" revision to Project2.docx
# Strategy: use Word's Find/Execute scoped to individual paragraphs (found via a
# unique anchor phrase) so that repeated short substrings ("0", "R", ...) can be
# targeted unambiguously, then toggle Bold (and occasionally re-typed text) on
# the precise sub-ranges the diff calls out.

$d = $word.ActiveDocument

function Get-ParaRange([string]$anchorText) {
    $a = $d.Content
    [void]$a.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $a.Paragraphs(1).Range
}

function Bold-SubRange([int]$pStart, [int]$pEnd, [string]$needle) {
    $r = $d.Range($pStart, $pEnd)
    [void]$r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Bold = 1
    return $r
}

# ---------------------------------------------------------------------------
# 1) "Can only move up, down, left and right." -- bold everything after "Can ",
#    keeping the existing italics on " up, down, left" / "right".
# ---------------------------------------------------------------------------
$para = Get-ParaRange("Can only move")
$pStart = $para.Start
$pEnd = $para.End

[void](Bold-SubRange $pStart $pEnd "only move")
[void](Bold-SubRange $pStart $pEnd " up, down, left")
[void](Bold-SubRange $pStart $pEnd " and ")
[void](Bold-SubRange $pStart $pEnd "right")

# ---------------------------------------------------------------------------
# 2) "... is 0~2147483647." -> "... is 0~2147483647(int)." with the whole
#    numeric range in bold (kept as several runs, mirroring how the original
#    "0" / "~" / number were already separate runs in the source document).
# ---------------------------------------------------------------------------
$para = Get-ParaRange("Note that the range of battery life")
$pStart = $para.Start
$pEnd = $para.End

$r = $d.Range($pStart, $pEnd)
[void]$r.Find.Execute("2147483647.", $true, $false, $false, $false, $false, $true, 1, $false, "2147483647(int).", 2)

# text length changed, re-acquire the paragraph before targeting sub-ranges
$para = Get-ParaRange("Note that the range of battery life")
$pStart = $para.Start
$pEnd = $para.End

[void](Bold-SubRange $pStart $pEnd "0")
[void](Bold-SubRange $pStart $pEnd "～")
[void](Bold-SubRange $pStart $pEnd "2147483647")
[void](Bold-SubRange $pStart $pEnd "(i")
[void](Bold-SubRange $pStart $pEnd "nt")
[void](Bold-SubRange $pStart $pEnd ")")

$r = $d.Range($pStart, $pEnd)
[void]$r.Find.Execute(").", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$periodRange = $d.Range($r.End - 1, $r.End)
$periodRange.Bold = 1

# ---------------------------------------------------------------------------
# 3) "Needs to return for recharge before running out battery." -- bold the
#    sentence, leave the rest of the paragraph untouched.
# ---------------------------------------------------------------------------
$para = Get-ParaRange("Needs to return for recharge")
$pStart = $para.Start
$pEnd = $para.End

[void](Bold-SubRange $pStart $pEnd "Needs to return for recharge before running out battery.")

# ---------------------------------------------------------------------------
# 4) & 5) "We assume the row index and column index all starts from 0 ..."
#    bold the first clause; later in the same paragraph bold "only one".
# ---------------------------------------------------------------------------
$para = Get-ParaRange("We assume the row index")
$pStart = $para.Start
$pEnd = $para.End

[void](Bold-SubRange $pStart $pEnd "We assume the row index and column index all starts from 0")
[void](Bold-SubRange $pStart $pEnd "only one")

# ---------------------------------------------------------------------------
# 6) "... should be no more than 1000*1000." -- bold "no more than 1000*1000".
# ---------------------------------------------------------------------------
$para = Get-ParaRange("Invalid test case gets")
$pStart = $para.Start
$pEnd = $para.End

[void](Bold-SubRange $pStart $pEnd "no more than 1000*1000")

# ---------------------------------------------------------------------------
# 7) Collapse the two runs split by the stray "_GoBack" bookmark back into a
#    single run (re-typing the same text forces Word to drop the bookmark).
# ---------------------------------------------------------------------------
$para = Get-ParaRange("Basic test (")
$pStart = $para.Start
$pEnd = $para.End

$r = $d.Range($pStart, $pEnd)
$phrase = "0%): If your test case can pass through TA" + [char]0x2019 + "s program, you receive "
[void]$r.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, $phrase, 2)

Write-Output "done"
